# Refresh the cryptocurrency table on Sheet1 with the latest scraped data
# (GitHub Actions run). Columns B/C identify the coin (name/link), D is the
# Price and E is Volume(1h); a few rows also swap which coin occupies them
# because the source ranking reordered.
#
# D and E are stored as plain TEXT in this sheet, not numbers -- prices use
# "." as a thousands separator (e.g. "34.670.65") and E keeps literal
# padding spaces around "+x.xx%". Assigning such strings straight to
# Range.Value lets Excel auto-detect some of them (e.g. "0.972") as real
# numbers, corrupting the data and changing the cell style. Set-TextValue
# forces a real text write (NumberFormat "@") and then restores the default
# "Normal" style, so only the cell value changes -- matching the diff, which
# shows no style/format changes at all.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") '34.670.65'
Set-TextValue $ws.Range("E2") '  +0.52%  '

# Row 3
Set-TextValue $ws.Range("D3") '1.829.07'
Set-TextValue $ws.Range("E3") '  +1.09%  '

# Row 4
Set-TextValue $ws.Range("E4") '  +0.13%  '

# Row 5
Set-TextValue $ws.Range("D5") '227.12'
Set-TextValue $ws.Range("E5") '  +0.58%  '

# Row 6
Set-TextValue $ws.Range("D6") '0.610'
Set-TextValue $ws.Range("E6") '  +1.89%  '

# Row 7
Set-TextValue $ws.Range("E7") '  +0.25%  '

# Row 8
Set-TextValue $ws.Range("D8") '44.02'
Set-TextValue $ws.Range("E8") '  +21.23%  '

# Row 9
Set-TextValue $ws.Range("E9") '  +1.98%  '

# Row 10
Set-TextValue $ws.Range("D10") '0.0686'
Set-TextValue $ws.Range("E10") '  +0.69%  '

# Row 11
Set-TextValue $ws.Range("E11") '  +3.80%  '

# Row 12
Set-TextValue $ws.Range("D12") '2.092.19'
Set-TextValue $ws.Range("E12") '  +1.02%  '

# Row 13
Set-TextValue $ws.Range("D13") '1.827.46'
Set-TextValue $ws.Range("E13") '  +0.89%  '

# Row 14
Set-TextValue $ws.Range("D14") '11.19'
Set-TextValue $ws.Range("E14") '  -1.15%  '

# Row 15
Set-TextValue $ws.Range("E15") '  +5.93%  '

# Row 16
Set-TextValue $ws.Range("D16") '0.650'
Set-TextValue $ws.Range("E16") '  +3.10%  '

# Row 17
Set-TextValue $ws.Range("D17") '34.652.38'
Set-TextValue $ws.Range("E17") '  +0.54%  '

# Row 18
Set-TextValue $ws.Range("D18") '68.31'
Set-TextValue $ws.Range("E18") '  -0.47%  '

# Row 19
Set-TextValue $ws.Range("D19") '242.75'
Set-TextValue $ws.Range("E19") '  -0.08%  '

# Row 20
Set-TextValue $ws.Range("D20") '0.0₃0784'
Set-TextValue $ws.Range("E20") '  +1.39%  '

# Row 21
Set-TextValue $ws.Range("D21") '12.14'
Set-TextValue $ws.Range("E21") '  +7.91%  '

# Row 22
Set-TextValue $ws.Range("D22") '4.73'
Set-TextValue $ws.Range("E22") '  +15.30%  '

# Row 23
Set-TextValue $ws.Range("E23") '  +0.19%  '

# Row 24
Set-TextValue $ws.Range("E24") '  -1.76%  '

# Row 25
Set-TextValue $ws.Range("D25") '171.33'
Set-TextValue $ws.Range("E25") '  +0.12%  '

# Row 26
Set-TextValue $ws.Range("D26") '7.95'
Set-TextValue $ws.Range("E26") '  +0.19%  '

# Row 27
Set-TextValue $ws.Range("D27") '17.73'
Set-TextValue $ws.Range("E27") '  +2.74%  '

# Row 28
Set-TextValue $ws.Range("E28") '  +0.71%  '

# Row 29
Set-TextValue $ws.Range("E29") '  +0.15%  '

# Row 30
Set-TextValue $ws.Range("B30") 'Filecoin'
Set-TextValue $ws.Range("C30") 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range("D30") '3.90'
Set-TextValue $ws.Range("E30") '  +2.23%  '

# Row 31
Set-TextValue $ws.Range("B31") 'PancakeSwap'
Set-TextValue $ws.Range("C31") 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue $ws.Range("D31") '1.26'
Set-TextValue $ws.Range("E31") '  +2.87%  '

# Row 32
Set-TextValue $ws.Range("D32") '3.97'
Set-TextValue $ws.Range("E32") '  +1.12%  '

# Row 33
Set-TextValue $ws.Range("D33") '0.0525'
Set-TextValue $ws.Range("E33") '  +1.46%  '

# Row 34
Set-TextValue $ws.Range("E34") '  +2.70%  '

# Row 35
Set-TextValue $ws.Range("D35") '90.04'
Set-TextValue $ws.Range("E35") '  +11.08%  '

# Row 36
Set-TextValue $ws.Range("D36") '0.666'
Set-TextValue $ws.Range("E36") '  +1.89%  '

# Row 37
Set-TextValue $ws.Range("D37") '1.332.92'
Set-TextValue $ws.Range("E37") '  -2.14%  '

# Row 38
Set-TextValue $ws.Range("D38") '2.44'
Set-TextValue $ws.Range("E38") '  +2.89%  '

# Row 39
Set-TextValue $ws.Range("D39") '15.28'
Set-TextValue $ws.Range("E39") '  +13.96%  '

# Row 40
Set-TextValue $ws.Range("E40") '  +0.20%  '

# Row 41
Set-TextValue $ws.Range("D41") '0.0192'
Set-TextValue $ws.Range("E41") '  +3.01%  '

# Row 42
Set-TextValue $ws.Range("B42") 'ARBITRUM'
Set-TextValue $ws.Range("C42") 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range("D42") '0.972'
Set-TextValue $ws.Range("E42") '  +3.72%  '

# Row 43
Set-TextValue $ws.Range("B43") 'WEMIXToken'
Set-TextValue $ws.Range("C43") 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws.Range("D43") '1.24'
Set-TextValue $ws.Range("E43") '  +6.17%  '

# Row 44
Set-TextValue $ws.Range("B44") 'MXToken'
Set-TextValue $ws.Range("C44") 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws.Range("D44") '2.81'
Set-TextValue $ws.Range("E44") '  +1.05%  '

# Row 45
Set-TextValue $ws.Range("B45") 'HuobiToken'
Set-TextValue $ws.Range("C45") 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue $ws.Range("D45") '2.43'
Set-TextValue $ws.Range("E45") '  +0.15%  '

# Row 46
Set-TextValue $ws.Range("D46") '0.0517'
Set-TextValue $ws.Range("E46") '  +3.56%  '

# Row 47
Set-TextValue $ws.Range("B47") 'FraxShare'
Set-TextValue $ws.Range("C47") 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range("D47") '5.97'
Set-TextValue $ws.Range("E47") '  +2.23%  '

# Row 48
Set-TextValue $ws.Range("B48") 'RocketPoolETH'
Set-TextValue $ws.Range("C48") 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextValue $ws.Range("D48") '1.990.28'
Set-TextValue $ws.Range("E48") '  +0.96%  '

# Row 49
Set-TextValue $ws.Range("E49") '  +0.32%  '

# Row 50
Set-TextValue $ws.Range("D50") '101.87'
Set-TextValue $ws.Range("E50") '  -0.86%  '

# Row 51
Set-TextValue $ws.Range("B51") 'Aptos'
Set-TextValue $ws.Range("C51") 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range("D51") '7.22'
Set-TextValue $ws.Range("E51") '  +3.62%  '
